# Update the "想去人数" (interested-attendee count) column F values across
# the 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets,
# per the gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibition) sheet: row -> new F value ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    4  = 88
    5  = 307
    7  = 139
    8  = 267
    9  = 23
    10 = 54
    12 = 127
    13 = 2762
    14 = 99
    15 = 36
    18 = 45
    19 = 554
    21 = 632
    22 = 189
    23 = 102
    25 = 31
    27 = 2274
    28 = 4818
    30 = 68
    31 = 471
    32 = 1249
    33 = 253
    34 = 2166
    37 = 70
    38 = 59
    39 = 142
    41 = 446
    42 = 753
    46 = 443
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# --- 演出 (Performance) sheet: row -> new F value ---
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    2 = 51
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Cells.Item($row, 6).Value = $sheet2Updates[$row]
}

# --- 全部类型 (All types) sheet: row -> new F value ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    4  = 88
    5  = 307
    7  = 139
    8  = 267
    9  = 23
    10 = 54
    12 = 127
    13 = 2762
    14 = 99
    15 = 36
    17 = 51
    19 = 45
    20 = 554
    22 = 632
    23 = 189
    24 = 102
    26 = 31
    28 = 2274
    29 = 4818
    31 = 68
    32 = 471
    33 = 1249
    34 = 253
    35 = 2166
    38 = 70
    39 = 59
    40 = 142
    42 = 446
    43 = 753
    47 = 443
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
